$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 4788.476
$ws.Range("I74").Value = 4853.6875
$ws.Range("J74").Value = 4579.8
$ws.Range("K74").Value = 4853.6875
$ws.Range("L74").Value = 4579.8
$ws.Range("M74").Value = -3917.6875
$ws.Range("N74").Value = -6451.8

# Row 76
$ws.Range("H76").Value = 3232.3333
$ws.Range("I76").Value = 2733.3333
$ws.Range("J76").Value = 3980.8333
$ws.Range("K76").Value = 2733.3333
$ws.Range("L76").Value = 3980.8333
$ws.Range("M76").Value = -2418.3333
$ws.Range("N76").Value = -4610.8333

# Row 77
$ws.Range("H77").Value = 4788.476
$ws.Range("I77").Value = 4853.6875
$ws.Range("J77").Value = 4579.8
$ws.Range("K77").Value = 24268.4375
$ws.Range("L77").Value = 22899
$ws.Range("M77").Value = -19588.4375
$ws.Range("N77").Value = -32259

# Row 79
$ws.Range("H79").Value = 3232.3333
$ws.Range("I79").Value = 2733.3333
$ws.Range("J79").Value = 3980.8333
$ws.Range("K79").Value = 2733.3333
$ws.Range("L79").Value = 3980.8333
$ws.Range("M79").Value = -1641.3333
$ws.Range("N79").Value = -6164.8333

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1248.0834
$ws.Range("I45").Value = 1291.6666
$ws.Range("K45").Value = 1291.6666
$ws.Range("M45").Value = -914.6666

# Row 63
$ws.Range("H63").Value = 2401.797
$ws.Range("I63").Value = 2398.6897
$ws.Range("J63").Value = 2418.182
$ws.Range("K63").Value = 2398.6897
$ws.Range("L63").Value = 2418.182
$ws.Range("M63").Value = -1712.6897
$ws.Range("N63").Value = -3790.182

# Row 66
$ws.Range("H66").Value = 2401.797
$ws.Range("I66").Value = 2398.6897
$ws.Range("J66").Value = 2418.182
$ws.Range("K66").Value = 11993.4485
$ws.Range("L66").Value = 12090.91
$ws.Range("M66").Value = -8561.448499999999
$ws.Range("N66").Value = -18954.91

# Row 74
$ws.Range("H74").Value = 21740544
$ws.Range("I74").Value = 25001256
$ws.Range("J74").Value = 2466.6667
$ws.Range("K74").Value = 25001256
$ws.Range("L74").Value = 2466.6667
$ws.Range("M74").Value = -25000382
$ws.Range("N74").Value = -4214.6667

# Row 77
$ws.Range("H77").Value = 21740544
$ws.Range("I77").Value = 25001256
$ws.Range("J77").Value = 2466.6667
$ws.Range("K77").Value = 125006280
$ws.Range("L77").Value = 12333.3335
$ws.Range("M77").Value = -125001912
$ws.Range("N77").Value = -21069.3335

# Row 105
$ws.Range("H105").Value = 20370
$ws.Range("J105").Value = 20370
$ws.Range("L105").Value = 20370
$ws.Range("N105").Value = -27358

# Row 122
$ws.Range("H122").Value = 1861.619
$ws.Range("I122").Value = 1499.5385
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 4498.6155
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -2048.6155
$ws.Range("N122").Value = -12250

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 366.55554
$ws.Range("I22").Value = 330.6154
$ws.Range("J22").Value = 460
$ws.Range("K22").Value = 330.6154
$ws.Range("L22").Value = 460
$ws.Range("M22").Value = -157.6154
$ws.Range("N22").Value = -806

# Row 105
$ws.Range("H105").Value = 2193.8572
$ws.Range("I105").Value = 1973.8462
$ws.Range("J105").Value = 2551.375
$ws.Range("K105").Value = 1973.8462
$ws.Range("L105").Value = 2551.375
$ws.Range("M105").Value = -226.8462
$ws.Range("N105").Value = -6045.375

# Row 110
$ws.Range("H110").Value = 21702
$ws.Range("J110").Value = 21702
$ws.Range("L110").Value = 21702
$ws.Range("N110").Value = -29882

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 19234350
$ws.Range("I31").Value = 41668750
$ws.Range("J31").Value = 4861.2856
$ws.Range("K31").Value = 41668750
$ws.Range("L31").Value = 4861.2856
$ws.Range("M31").Value = -41668455
$ws.Range("N31").Value = -5451.2856

# Row 34
$ws.Range("H34").Value = 19234350
$ws.Range("I34").Value = 41668750
$ws.Range("J34").Value = 4861.2856
$ws.Range("K34").Value = 41668750
$ws.Range("L34").Value = 4861.2856
$ws.Range("M34").Value = -41668548
$ws.Range("N34").Value = -5265.2856

# Row 62
$ws.Range("H62").Value = 92861.73
$ws.Range("I62").Value = 101647.9
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 101647.9
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -101023.9
$ws.Range("N62").Value = -6248

# Row 65
$ws.Range("H65").Value = 92861.73
$ws.Range("I65").Value = 101647.9
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 508239.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -505119.5
$ws.Range("N65").Value = -31240

# Row 118
$ws.Range("H118").Value = 43401.816
$ws.Range("J118").Value = 43401.816
$ws.Range("L118").Value = 43401.816
$ws.Range("N118").Value = -46715.816

# Row 134
$ws.Range("H134").Value = 1251.5172
$ws.Range("I134").Value = 896.7619
$ws.Range("J134").Value = 2182.75
$ws.Range("K134").Value = 2690.2857
$ws.Range("L134").Value = 6548.25
$ws.Range("M134").Value = -155.2856999999999
$ws.Range("N134").Value = -11618.25

$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 4029.7058
$ws.Range("J51").Value = 4961.923
$ws.Range("L51").Value = 14885.769
$ws.Range("N51").Value = -15805.769

# Row 110
$ws.Range("H110").Value = 4750
$ws.Range("I110").Value = 4750
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 14250
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -10160
$ws.Range("N110").ClearContents()

# Row 131
$ws.Range("H131").Value = 842.78845
$ws.Range("I131").Value = 273.2
$ws.Range("J131").Value = 978.4048
$ws.Range("K131").Value = 819.5999999999999
$ws.Range("L131").Value = 2935.2144
$ws.Range("M131").Value = 4220.4
$ws.Range("N131").Value = -13015.2144

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6476.75
$ws.Range("I70").Value = 7978.8
$ws.Range("J70").Value = 3973.3333
$ws.Range("K70").Value = 7978.8
$ws.Range("L70").Value = 3973.3333
$ws.Range("M70").Value = -7708.8
$ws.Range("N70").Value = -4513.3333

# Row 73
$ws.Range("H73").Value = 6476.75
$ws.Range("I73").Value = 7978.8
$ws.Range("J73").Value = 3973.3333
$ws.Range("K73").Value = 7978.8
$ws.Range("L73").Value = 3973.3333
$ws.Range("M73").Value = -7042.8
$ws.Range("N73").Value = -5845.3333

# Row 80
$ws.Range("H80").Value = 43725
$ws.Range("I80").Value = 2035.5
$ws.Range("J80").Value = 73503.21000000001
$ws.Range("K80").Value = 2035.5
$ws.Range("L80").Value = 73503.21000000001
$ws.Range("M80").Value = -1037.5
$ws.Range("N80").Value = -75499.21000000001

# Row 83
$ws.Range("H83").Value = 43725
$ws.Range("I83").Value = 2035.5
$ws.Range("J83").Value = 73503.21000000001
$ws.Range("K83").Value = 10177.5
$ws.Range("L83").Value = 367516.05
$ws.Range("M83").Value = -5185.5
$ws.Range("N83").Value = -377500.05

# Row 104
$ws.Range("H104").Value = 36835.5
$ws.Range("J104").Value = 36835.5
$ws.Range("L104").Value = 36835.5
$ws.Range("N104").Value = -43823.5

# Row 122
$ws.Range("H122").Value = 2322.5386
$ws.Range("I122").Value = 2540.5715
$ws.Range("J122").Value = 2068.1667
$ws.Range("K122").Value = 7621.7145
$ws.Range("L122").Value = 6204.500100000001
$ws.Range("M122").Value = -5171.7145
$ws.Range("N122").Value = -11104.5001

# Row 123
$ws.Range("H123").Value = 20199.334
$ws.Range("J123").Value = 20199.334
$ws.Range("L123").Value = 20199.334
$ws.Range("N123").Value = -25099.334

$ws = $wb.Worksheets.Item("LTW")
# Row 29
$ws.Range("H29").Value = 7503.2
$ws.Range("I29").Value = 1016
$ws.Range("J29").Value = 9125
$ws.Range("K29").Value = 1016
$ws.Range("L29").Value = 9125
$ws.Range("M29").Value = -721
$ws.Range("N29").Value = -9715

# Row 110
$ws.Range("H110").Value = 16186
$ws.Range("I110").Value = 9900
$ws.Range("J110").Value = 18281.334
$ws.Range("K110").Value = 9900
$ws.Range("L110").Value = 18281.334
$ws.Range("M110").Value = -5810
$ws.Range("N110").Value = -26461.334
